$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine the used range of the sheet.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row  # xlUp

# Swap the "grade" column (B) with the "name" column (C) -- the sheet is
# being restructured so column B becomes "course_type" (Traditional/Online)
# and column C becomes "grade" (A/B/C/D/F). The underlying per-row data
# doesn't change, only which column it lives in.
$rangeB = $ws.Range("B2:B$lastRow")
$rangeC = $ws.Range("C2:C$lastRow")

$valuesB = $rangeB.Value2
$valuesC = $rangeC.Value2

$rangeB.Value2 = $valuesC
$rangeC.Value2 = $valuesB

# Update the header row to reflect the new column meanings.
$ws.Range("B1").Value2 = "course_type"
$ws.Range("C1").Value2 = "grade"

# Column B now holds the longer "Traditional"/"Online" text instead of a
# single-letter grade, so widen it (best-fit) to match.
$ws.Columns.Item(2).ColumnWidth = 10.25
